# Apply filter based on sequence coverage and unique peptides:
# source/target link table is re-filtered and re-ordered (Mrpl4 rows first,
# then Cbp3 rows), fold-change values are refreshed, and one new row is
# added (49 -> 50 data+header rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 49,3
$data[0,0] = "Mrpl4"
$data[0,1] = "Mgr1"
$data[0,2] = 2.242187785
$data[1,0] = "Mrpl4"
$data[1,1] = "Yfr045w"
$data[1,2] = 4.476406151
$data[2,0] = "Mrpl4"
$data[2,1] = "Cyt1"
$data[2,2] = 2.627029822
$data[3,0] = "Mrpl4"
$data[3,1] = "Cbp6"
$data[3,2] = 6.073919856
$data[4,0] = "Mrpl4"
$data[4,1] = "Pet111"
$data[4,2] = 1.779152136
$data[5,0] = "Mrpl4"
$data[5,1] = "Sal1"
$data[5,2] = 2.500020235
$data[6,0] = "Mrpl4"
$data[6,1] = "Cbp3"
$data[6,2] = 4.525702317
$data[7,0] = "Mrpl4"
$data[7,1] = "Mrp20"
$data[7,2] = 4.153036304
$data[8,0] = "Mrpl4"
$data[8,1] = "Nat2"
$data[8,2] = 4.078817589
$data[9,0] = "Mrpl4"
$data[9,1] = "Mba1"
$data[9,2] = 2.061652968
$data[10,0] = "Mrpl4"
$data[10,1] = "Mcx1"
$data[10,2] = 1.686755299
$data[11,0] = "Mrpl4"
$data[11,1] = "Cox15"
$data[11,2] = 2.378184154
$data[12,0] = "Mrpl4"
$data[12,1] = "Mmf1"
$data[12,2] = 1.749621579
$data[13,0] = "Mrpl4"
$data[13,1] = "Nde1"
$data[13,2] = 2.154364779
$data[14,0] = "Mrpl4"
$data[14,1] = "Coa1"
$data[14,2] = 5.304607176
$data[15,0] = "Mrpl4"
$data[15,1] = "Tes1"
$data[15,2] = 3.103354815
$data[16,0] = "Mrpl4"
$data[16,1] = "Ydl183c"
$data[16,2] = 2.139772469
$data[17,0] = "Mrpl4"
$data[17,1] = "Ynr040w"
$data[17,2] = 2.644640891
$data[18,0] = "Mrpl4"
$data[18,1] = "Tim50"
$data[18,2] = 3.091846855
$data[19,0] = "Mrpl4"
$data[19,1] = "Ydl027c"
$data[19,2] = 2.620800279
$data[20,0] = "Mrpl4"
$data[20,1] = "Rdl2"
$data[20,2] = 2.426825355
$data[21,0] = "Mrpl4"
$data[21,1] = "Ypl168w"
$data[21,2] = 6.315955929
$data[22,0] = "Cbp3"
$data[22,1] = "Mgr1"
$data[22,2] = 3.037622469
$data[23,0] = "Cbp3"
$data[23,1] = "Yfr045w"
$data[23,2] = 4.291550761
$data[24,0] = "Cbp3"
$data[24,1] = "Aim39"
$data[24,2] = 2.005386778
$data[25,0] = "Cbp3"
$data[25,1] = "Pgk1"
$data[25,2] = 2.687371265
$data[26,0] = "Cbp3"
$data[26,1] = "Cyt1"
$data[26,2] = 4.509850964
$data[27,0] = "Cbp3"
$data[27,1] = "Cbp6"
$data[27,2] = 4.358683278
$data[28,0] = "Cbp3"
$data[28,1] = "Lpd1"
$data[28,2] = 1.575035635
$data[29,0] = "Cbp3"
$data[29,1] = "Sal1"
$data[29,2] = 2.591443261
$data[30,0] = "Cbp3"
$data[30,1] = "Cbp3"
$data[30,2] = 8.224277944000001
$data[31,0] = "Cbp3"
$data[31,1] = "Yme1"
$data[31,2] = 2.400876254
$data[32,0] = "Cbp3"
$data[32,1] = "Nat2"
$data[32,2] = 3.526374914
$data[33,0] = "Cbp3"
$data[33,1] = "Mba1"
$data[33,2] = 1.888152295
$data[34,0] = "Cbp3"
$data[34,1] = "Mcx1"
$data[34,2] = 1.813332533
$data[35,0] = "Cbp3"
$data[35,1] = "Ylf2"
$data[35,2] = 2.999642688
$data[36,0] = "Cbp3"
$data[36,1] = "Cox15"
$data[36,2] = 2.279966902
$data[37,0] = "Cbp3"
$data[37,1] = "Mmf1"
$data[37,2] = 1.784193007
$data[38,0] = "Cbp3"
$data[38,1] = "Nde1"
$data[38,2] = 3.061170275
$data[39,0] = "Cbp3"
$data[39,1] = "Coa1"
$data[39,2] = 5.19063741
$data[40,0] = "Cbp3"
$data[40,1] = "Tes1"
$data[40,2] = 2.913472794
$data[41,0] = "Cbp3"
$data[41,1] = "Ydl183c"
$data[41,2] = 2.331665432
$data[42,0] = "Cbp3"
$data[42,1] = "Ynr040w"
$data[42,2] = 1.629489825
$data[43,0] = "Cbp3"
$data[43,1] = "Tim50"
$data[43,2] = 3.938119285
$data[44,0] = "Cbp3"
$data[44,1] = "Ydl027c"
$data[44,2] = 2.98248517
$data[45,0] = "Cbp3"
$data[45,1] = "Fmp25"
$data[45,2] = 3.94528311
$data[46,0] = "Cbp3"
$data[46,1] = "Mdm38"
$data[46,2] = 1.837488681
$data[47,0] = "Cbp3"
$data[47,1] = "Rdl2"
$data[47,2] = 1.798418196
$data[48,0] = "Cbp3"
$data[48,1] = "Ypl168w"
$data[48,2] = 5.757088287

$ws.Range("A2:C50").Value = $data

Write-Output "done"
